$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 32.51511900000001
$ws.Range("H2").Value = 97.54535700000001
$ws.Range("I2").Value = 0.218203973858649
$ws.Range("J2").Value = 0.2182039738586489
$ws.Range("M2").Value = 0.2595153333333333
$ws.Range("N2").Value = 0.778546
$ws.Range("O2").Value = 0.008853914448786948
$ws.Range("P2").Value = 0.008853914448786946
$ws.Range("Q2").Value = 8.438171945658
$ws.Range("R2").Value = 75.943547510922
$ws.Range("S2").Value = 0.001931959316929821
$ws.Range("T2").Value = 0.001931959316929821
$ws.Range("G3").Value = 32.51511900000001
$ws.Range("H3").Value = 97.54535700000001
$ws.Range("I3").Value = 0.218203973858649
$ws.Range("J3").Value = 0.2182039738586489
$ws.Range("O3").Value = 0.80858400485482
$ws.Range("P3").Value = 0.80858400485482
$ws.Range("Q3").Value = 770.616308180901
$ws.Range("R3").Value = 6935.546773628109
$ws.Range("S3").Value = 0.1764362430578628
$ws.Range("T3").Value = 0.1764362430578628
$ws.Range("G4").Value = 32.51511900000001
$ws.Range("H4").Value = 97.54535700000001
$ws.Range("I4").Value = 0.218203973858649
$ws.Range("J4").Value = 0.2182039738586489
$ws.Range("M4").Value = 5.351040999999999
$ws.Range("N4").Value = 16.053123
$ws.Range("O4").Value = 0.1825620806963931
$ws.Range("P4").Value = 0.1825620806963931
$ws.Range("Q4").Value = 173.989734888879
$ws.Range("R4").Value = 1565.907613999911
$ws.Range("S4").Value = 0.03983577148385632
$ws.Range("T4").Value = 0.03983577148385631
$ws.Range("I5").Value = 0.6017421411306194
$ws.Range("J5").Value = 0.6017421411306194
$ws.Range("M5").Value = 0.2595153333333333
$ws.Range("N5").Value = 0.778546
$ws.Range("O5").Value = 0.008853914448786948
$ws.Range("P5").Value = 0.008853914448786946
$ws.Range("Q5").Value = 23.269987086018
$ws.Range("R5").Value = 209.429883774162
$ws.Range("S5").Value = 0.005327773437800386
$ws.Range("T5").Value = 0.005327773437800385
$ws.Range("I6").Value = 0.6017421411306194
$ws.Range("J6").Value = 0.6017421411306194
$ws.Range("O6").Value = 0.80858400485482
$ws.Range("P6").Value = 0.80858400485482
$ws.Range("S6").Value = 0.4865590703653105
$ws.Range("T6").Value = 0.4865590703653105
$ws.Range("I7").Value = 0.6017421411306194
$ws.Range("J7").Value = 0.6017421411306194
$ws.Range("M7").Value = 5.351040999999999
$ws.Range("N7").Value = 16.053123
$ws.Range("O7").Value = 0.1825620806963931
$ws.Range("P7").Value = 0.1825620806963931
$ws.Range("Q7").Value = 479.812323100059
$ws.Range("R7").Value = 4318.310907900531
$ws.Range("S7").Value = 0.1098552973275085
$ws.Range("T7").Value = 0.1098552973275085
$ws.Range("G8").Value = 26.830279
$ws.Range("H8").Value = 80.490837
$ws.Range("I8").Value = 0.1800538850107317
$ws.Range("J8").Value = 0.1800538850107317
$ws.Range("M8").Value = 0.2595153333333333
$ws.Range("N8").Value = 0.778546
$ws.Range("O8").Value = 0.008853914448786948
$ws.Range("P8").Value = 0.008853914448786946
$ws.Range("Q8").Value = 6.962868798111333
$ws.Range("R8").Value = 62.66581918300199
$ws.Range("S8").Value = 0.001594181694056741
$ws.Range("T8").Value = 0.00159418169405674
$ws.Range("G9").Value = 26.830279
$ws.Range("H9").Value = 80.490837
$ws.Range("I9").Value = 0.1800538850107317
$ws.Range("J9").Value = 0.1800538850107317
$ws.Range("O9").Value = 0.80858400485482
$ws.Range("P9").Value = 0.80858400485482
$ws.Range("Q9").Value = 635.8842036052076
$ws.Range("R9").Value = 5722.957832446868
$ws.Range("S9").Value = 0.1455886914316467
$ws.Range("T9").Value = 0.1455886914316467
$ws.Range("G10").Value = 26.830279
$ws.Range("H10").Value = 80.490837
$ws.Range("I10").Value = 0.1800538850107317
$ws.Range("J10").Value = 0.1800538850107317
$ws.Range("M10").Value = 5.351040999999999
$ws.Range("N10").Value = 16.053123
$ws.Range("O10").Value = 0.1825620806963931
$ws.Range("P10").Value = 0.1825620806963931
$ws.Range("Q10").Value = 143.569922970439
$ws.Range("R10").Value = 1292.129306733951
$ws.Range("S10").Value = 0.03287101188502829
$ws.Range("T10").Value = 0.03287101188502827
